$d = $word.ActiveDocument

function Replace-Text {
    param($findText, $replaceText)
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "WARNING: not found -> $findText"
    }
}

# --- Main paragraph: grammar fixes, added detail (Stokes/anti-Stokes, sub-diffraction-limit, etc.) ---
Replace-Text "od access the absolute temperature of go" "od accesses the absolute temperature of go"
Replace-Text "ld nanoparticles by measuring the photolum" "ld nanoparticles by measuring their Stokes and anti-Stokes photolum"
Replace-Text "inescence emission spectra. Therefo" "inescence spectra. Therefo"
Replace-Text "re, it capitalizes the unique optical properties of gold n" "re, it capitalizes on the unique optical properties of gold n"
Replace-Text "anoparticles such as stable photoluminesnce emission and their nanometr" "anoparticles such as stable photoluminescence and their nanometr"
Replace-Text "ic size to probe subdifrac" "ic size to probe sub-diffrac"
Replace-Text "tion volumes. More important" "tion-limit volumes. More important"
Replace-Text "ly, it has the advantage of avoiding temper" "ly, it has the advantages of avoiding temper"
Replace-Text "ature calibrations, being non-invasive and easy to implem" "ature calibrations, of being non-invasive and easy to implem"
Replace-Text "ent in a regular microscope with spectral capabilities. As a proof of pr" "ent in a regular microscope with spectrometric capabilities. As a proof of pr"
Replace-Text "inciple, we used the method to optically access the te" "inciple, we used the method to access the te"
Replace-Text "mperature of resonantly-illuminated gold nan" "mperature of the nan"
Replace-Text "orods and the temp" "orods themselves and the temp"
Replace-Text "erature of their surroun" "erature of the surroun"
Replace-Text "ding media with an ac" "ding medium with an ac"

# --- Trim the closing sentence of the third paragraph ---
Replace-Text " and that this journal will be the appropriate dissemination vehicle for it." "."

# --- Move the _GoBack bookmark to the last edited location ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$anchor = $d.Content
[void]$anchor.Find.Execute("broad audience of Nano Letters")
$insertPoint = $d.Range($anchor.End, $anchor.End)
[void]$d.Bookmarks.Add("_GoBack", $insertPoint)

# --- Normalize the "Thank you" run split (merge identical-format runs) ---
Replace-Text "Thank you for your consideration of our work." "Thank you for your consideration of our work."
